$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update country names (shared-string reorder effect) in column A ---
$ws.Range("A29").Value = "Ecuador"
$ws.Range("A30").Value = "Bolivia"
$ws.Range("A44").Value = "Guatemala"
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("A46").Value = "Paises Bajos"
$ws.Range("A63").Value = "Etiopia"
$ws.Range("A64").Value = "Azerbaiyan"
$ws.Range("A87").Value = "Libano"
$ws.Range("A88").Value = "Zambia"
$ws.Range("A123").Value = "Mozambique"
$ws.Range("A124").Value = "Eslovaquia"

# --- Update "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 19:38"

# --- Update numeric data cells (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5717825
$ws.Range("C4").Value = 16894
$ws.Range("D4").Value = 3064631
$ws.Range("E4").Value = 2476465
$ws.Range("G4").Value = 392
$ws.Range("H4").Value = 176729
$ws.Range("B5").Value = 3470517
$ws.Range("C5").Value = 10104
$ws.Range("E5").Value = 743820
$ws.Range("G5").Value = 254
$ws.Range("H5").Value = 111443
$ws.Range("B6").Value = 2903676
$ws.Range("C6").Value = 67854
$ws.Range("D6").Value = 2157187
$ws.Range("E6").Value = 691526
$ws.Range("G6").Value = 969
$ws.Range("H6").Value = 54963
$ws.Range("B22").Value = 230183
$ws.Range("C22").Value = 483
$ws.Range("E22").Value = 16065
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 9318
$ws.Range("B23").Value = 229814
$ws.Range("C23").Value = 4771
$ws.Range("E23").Value = 115269
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 30480
$ws.Range("B29").Value = 105508
$ws.Range("C29").Value = 1033
$ws.Range("D29").Value = 87660
$ws.Range("E29").Value = 11648
$ws.Range("G29").Value = 54
$ws.Range("H29").Value = 6200
$ws.Range("B30").Value = 105050
$ws.Range("C30").Value = 2031
$ws.Range("D30").Value = 39965
$ws.Range("E30").Value = 60852
$ws.Range("G30").Value = 61
$ws.Range("H30").Value = 4233
$ws.Range("B32").Value = 99201
$ws.Range("C32").Value = 1232
$ws.Range("D32").Value = 74536
$ws.Range("E32").Value = 23870
$ws.Range("G32").Value = 14
$ws.Range("H32").Value = 795
$ws.Range("B43").Value = 69950
$ws.Range("C43").Value = 149
$ws.Range("D43").Value = 67929
$ws.Range("E43").Value = 1394
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 627
$ws.Range("B44").Value = 65983
$ws.Range("C44").Value = 1102
$ws.Range("D44").Value = 54351
$ws.Range("E44").Value = 9126
$ws.Range("G44").Value = 39
$ws.Range("H44").Value = 2506
$ws.Range("B45").Value = 65802
$ws.Range("C45").Value = 461
$ws.Range("D45").Value = 58153
$ws.Range("E45").Value = 7280
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 369
$ws.Range("B46").Value = 65054
$ws.Range("C46").Value = 529
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("G46").Value = 10
$ws.Range("H46").Value = 6191
$ws.Range("B54").Value = 47638
$ws.Range("C54").Value = 1325
$ws.Range("D54").Value = 32806
$ws.Range("E54").Value = 14057
$ws.Range("G54").Value = 32
$ws.Range("H54").Value = 775
$ws.Range("B58").Value = 40258
$ws.Range("C58").Value = 411
$ws.Range("D58").Value = 28281
$ws.Range("E58").Value = 10566
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 1411
$ws.Range("B63").Value = 35836
$ws.Range("C63").Value = 1778
$ws.Range("D63").Value = 13536
$ws.Range("E63").Value = 21680
$ws.Range("G63").Value = 20
$ws.Range("H63").Value = 620
$ws.Range("B64").Value = 34759
$ws.Range("C64").Value = 139
$ws.Range("D64").Value = 32511
$ws.Range("E64").Value = 1738
$ws.Range("H64").Value = 510
$ws.Range("B87").Value = 10952
$ws.Range("C87").Value = 605
$ws.Range("D87").Value = 3040
$ws.Range("E87").Value = 7799
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 113
$ws.Range("B88").Value = 10372
$ws.Range("C88").Value = 154
$ws.Range("D88").Value = 9126
$ws.Range("E88").Value = 972
$ws.Range("G88").Value = 5
$ws.Range("H88").Value = 274
$ws.Range("B94").Value = 8743
$ws.Range("C94").Value = 32
$ws.Range("D94").Value = 8189
$ws.Range("E94").Value = 500
$ws.Range("B104").Value = 6370
$ws.Range("C104").Value = 145
$ws.Range("D104").Value = 3915
$ws.Range("E104").Value = 2431
$ws.Range("B123").Value = 3115
$ws.Range("C123").Value = 70
$ws.Range("D123").Value = 1380
$ws.Range("E123").Value = 1715
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 20
$ws.Range("B124").Value = 3102
$ws.Range("C124").Value = 80
$ws.Range("D124").Value = 2014
$ws.Range("E124").Value = 1055
$ws.Range("H124").Value = 33
$ws.Range("B134").Value = 2149
$ws.Range("C134").Value = 32
$ws.Range("E134").Value = 1101
$ws.Range("B138").Value = 1969
$ws.Range("C138").Value = 8
$ws.Range("D138").Value = 1536
$ws.Range("E138").Value = 364
$ws.Range("B140").Value = 1899
$ws.Range("C140").Value = 7
$ws.Range("D140").Value = 1058
$ws.Range("E140").Value = 300
$ws.Range("G140").Value = 2
$ws.Range("H140").Value = 541
$ws.Range("B181").Value = 306
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 274
$ws.Range("E181").Value = 32
$ws.Range("B184").Value = 229
$ws.Range("C184").Value = 6
$ws.Range("D184").Value = 200
$ws.Range("E184").Value = 29
